# Weekly update: insert a new price record as row 188, pushing the
# existing rows 188-211 down to 189-212 (dimension grows to A1:R212).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 188 (shifts 188..211 -> 189..212).
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with the new weekly record.
$ws.Range("A188").Value = 7
$ws.Range("B188").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C188").Value = "Ñuble"
$ws.Range("D188").Value = 44491
$ws.Range("E188").Value = 16
$ws.Range("F188").Value = 100114001
$ws.Range("G188").Value = "Papa"
$ws.Range("H188").Value = "Patagonia"
$ws.Range("I188").Value = "1a (guarda lavada)"
$ws.Range("J188").Value = 180
$ws.Range("K188").Value = 9500
$ws.Range("L188").Value = 10000
$ws.Range("M188").Value = 9750
$ws.Range("N188").Value = "$/malla 25 kilos"
$ws.Range("O188").Value = "Región de Los Lagos"
$ws.Range("P188").Value = 390
$ws.Range("Q188").Value = 25
$ws.Range("R188").Value = "Hortaliza"
